# --- Add "Sheet2" (face-recognition attendance register) after Sheet1, make it the active sheet ---
$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# --- Column widths (matches the "best fit" widths of the attendance table) ---
$ws2.Columns.Item(1).ColumnWidth = 7.833333333333333
$ws2.Columns.Item(2).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(3).ColumnWidth = 7.833333333333333
$ws2.Columns.Item(4).ColumnWidth = 9.333333333333332
$ws2.Columns.Item(5).ColumnWidth = 3.333333333333333
$ws2.Columns.Item(6).ColumnWidth = 5.0
$ws2.Columns.Item(7).ColumnWidth = 3.1666666666666665
$ws2.Columns.Item(8).ColumnWidth = 4.0
$ws2.Columns.Item(9).ColumnWidth = 4.333333333333333
$ws2.Columns.Item(10).ColumnWidth = 2.1666666666666665
$ws2.Columns.Item(11).ColumnWidth = 2.6666666666666665
$ws2.Columns.Item(12).ColumnWidth = 3.333333333333333
$ws2.Columns.Item(13).ColumnWidth = 5.0
$ws2.Columns.Item(14).ColumnWidth = 3.1666666666666665
$ws2.Columns.Item(15).ColumnWidth = 4.0
$ws2.Columns.Item(16).ColumnWidth = 4.333333333333333
$ws2.Columns.Item(17).ColumnWidth = 2.1666666666666665
$ws2.Columns.Item(18).ColumnWidth = 2.6666666666666665
$ws2.Columns.Item(19).ColumnWidth = 3.333333333333333

# --- Header rows: gr_no / fname / lname / surname, day numbers 1-15, weekday labels ---
$ws2.Range("A2").Value = "gr_no"
$ws2.Range("B2").Value = "fname"
$ws2.Range("C2").Value = "lname"
$ws2.Range("D2").Value = "surname"
$ws2.Range("E2").Value = 1
$ws2.Range("F2").Value = 2
$ws2.Range("G2").Value = 3
$ws2.Range("H2").Value = 4
$ws2.Range("I2").Value = 5
$ws2.Range("J2").Value = 6
$ws2.Range("K2").Value = 7
$ws2.Range("L2").Value = 8
$ws2.Range("M2").Value = 9
$ws2.Range("N2").Value = 10
$ws2.Range("O2").Value = 11
$ws2.Range("P2").Value = 12
$ws2.Range("Q2").Value = 13
$ws2.Range("R2").Value = 14
$ws2.Range("S2").Value = 15

$ws2.Range("E3").Value = "sun"
$ws2.Range("F3").Value = "mon  "
$ws2.Range("G3").Value = "tue"
$ws2.Range("H3").Value = "wed"
$ws2.Range("I3").Value = "thue"
$ws2.Range("J3").Value = "fri"
$ws2.Range("K3").Value = "sat"
$ws2.Range("L3").Value = "sun"
$ws2.Range("M3").Value = "mon  "
$ws2.Range("N3").Value = "tue"
$ws2.Range("O3").Value = "wed"
$ws2.Range("P3").Value = "thue"
$ws2.Range("Q3").Value = "fri"
$ws2.Range("R3").Value = "sat"
$ws2.Range("S3").Value = "sun"

$ws2.Range("A5").Value = "2025-001"
$ws2.Range("B5").Value = "Catherine"
$ws2.Range("C5").Value = "Marshall"
$ws2.Range("D5").Value = "Ware"
$ws2.Range("E5").Value = "P"
$ws2.Range("F5").Value = "P"
$ws2.Range("G5").Value = "P"
$ws2.Range("H5").Value = "P"
$ws2.Range("I5").Value = "P"
$ws2.Range("J5").Value = "P"
$ws2.Range("K5").Value = "P"
$ws2.Range("L5").Value = "P"
$ws2.Range("M5").Value = "P"
$ws2.Range("N5").Value = "P"
$ws2.Range("O5").Value = "P"
$ws2.Range("P5").Value = "P"
$ws2.Range("Q5").Value = "P"
$ws2.Range("R5").Value = "P"
$ws2.Range("S5").Value = "P"

$ws2.Range("A6").Value = "2025-002"
$ws2.Range("B6").Value = "Dawn"
$ws2.Range("C6").Value = "Smith"
$ws2.Range("D6").Value = "Valencia"
$ws2.Range("E6").Value = "A"
$ws2.Range("F6").Value = "A"
$ws2.Range("G6").Value = "A"
$ws2.Range("H6").Value = "A"
$ws2.Range("I6").Value = "A"
$ws2.Range("J6").Value = "A"
$ws2.Range("K6").Value = "A"
$ws2.Range("L6").Value = "A"
$ws2.Range("M6").Value = "A"
$ws2.Range("N6").Value = "A"
$ws2.Range("O6").Value = "A"
$ws2.Range("P6").Value = "A"
$ws2.Range("Q6").Value = "A"
$ws2.Range("R6").Value = "A"
$ws2.Range("S6").Value = "A"

$ws2.Range("A7").Value = "2025-003"
$ws2.Range("B7").Value = "Casey"
$ws2.Range("C7").Value = "Serrano"
$ws2.Range("D7").Value = "Rogers"

$ws2.Range("A8").Value = "2025-004"
$ws2.Range("B8").Value = "Eduardo"
$ws2.Range("C8").Value = "Bowers"
$ws2.Range("D8").Value = "Howard"

$ws2.Range("A9").Value = "2025-005"
$ws2.Range("B9").Value = "Michele"
$ws2.Range("C9").Value = "Warren"
$ws2.Range("D9").Value = "Shields"

$ws2.Range("A10").Value = "2025-006"
$ws2.Range("B10").Value = "Zachary"
$ws2.Range("C10").Value = "Riley"
$ws2.Range("D10").Value = "Campbell"

$ws2.Range("A11").Value = "2025-007"
$ws2.Range("B11").Value = "Amber"
$ws2.Range("C11").Value = "Marquez"
$ws2.Range("D11").Value = "Ramirez"

$ws2.Range("A12").Value = "2025-008"
$ws2.Range("B12").Value = "Kimberly"
$ws2.Range("C12").Value = "Shannon"
$ws2.Range("D12").Value = "Vazquez"

$ws2.Range("A13").Value = "2025-009"
$ws2.Range("B13").Value = "Jennifer"
$ws2.Range("C13").Value = "Cross"
$ws2.Range("D13").Value = "Mitchell"

$ws2.Range("A14").Value = "2025-010"
$ws2.Range("B14").Value = "Emily"
$ws2.Range("C14").Value = "Ryan"
$ws2.Range("D14").Value = "Mills"

$ws2.Range("A15").Value = "2025-011"
$ws2.Range("B15").Value = "Amanda"
$ws2.Range("C15").Value = "Johnson"
$ws2.Range("D15").Value = "Keller"

$ws2.Range("A16").Value = "2025-012"
$ws2.Range("B16").Value = "Anthony"
$ws2.Range("C16").Value = "Bradley"
$ws2.Range("D16").Value = "Crawford"

$ws2.Range("A17").Value = "2025-013"
$ws2.Range("B17").Value = "Tonya"
$ws2.Range("C17").Value = "Peters"
$ws2.Range("D17").Value = "Fuller"

$ws2.Range("A18").Value = "2025-014"
$ws2.Range("B18").Value = "Tom"
$ws2.Range("C18").Value = "Davis"
$ws2.Range("D18").Value = "Hammond"

$ws2.Range("A19").Value = "2025-015"
$ws2.Range("B19").Value = "Savannah"
$ws2.Range("C19").Value = "Williams"
$ws2.Range("D19").Value = "Jackson"

$ws2.Range("A20").Value = "2025-016"
$ws2.Range("B20").Value = "Donna"
$ws2.Range("C20").Value = "Zuniga"
$ws2.Range("D20").Value = "Gilbert"

$ws2.Range("A21").Value = "2025-017"
$ws2.Range("B21").Value = "Christopher"
$ws2.Range("C21").Value = "Roberts"
$ws2.Range("D21").Value = "Diaz"

$ws2.Range("A22").Value = "2025-018"
$ws2.Range("B22").Value = "Julie"
$ws2.Range("C22").Value = "Jensen"
$ws2.Range("D22").Value = "Hamilton"

$ws2.Range("A23").Value = "2025-019"
$ws2.Range("B23").Value = "Laura"
$ws2.Range("C23").Value = "Thomas"
$ws2.Range("D23").Value = "Gray"

$ws2.Range("A24").Value = "2025-020"
$ws2.Range("B24").Value = "Keith"
$ws2.Range("C24").Value = "Taylor"
$ws2.Range("D24").Value = "Lopez"

$ws2.Range("A25").Value = "2021-013"
$ws2.Range("B25").Value = "Bhavani "
$ws2.Range("C25").Value = "J"
$ws2.Range("D25").Value = "Gohil"

$ws2.Range("A26").Value = "2019-008"
$ws2.Range("B26").Value = "yash"
$ws2.Range("C26").Value = "R"
$ws2.Range("D26").Value = "prajapati"

# --- Selection left where the author left it when saving ---
$ws2.Range("U9").Select()

Write-Output "Sheet2 (attendance) created"